$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.613.10"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.893.57"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.07"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2939"
$ws.Range("E8").Value = "  +2.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06719"
$ws.Range("E9").Value = "  +0.99%  "

$ws.Range("D10").Value = "1.912.44"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.13"
$ws.Range("E11").Value = "  +2.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07354"
$ws.Range("E12").Value = "  +1.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.163"
$ws.Range("E13").Value = "  +3.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.17"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6704"
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").Value = "30.546.37"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007868"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("E18").Value = "  +3.64%  "

$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.146.04"
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.292"
$ws.Range("E21").Value = "  +11.87%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "191.05"
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.183"
$ws.Range("E24").Value = "  +2.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.494"
$ws.Range("E25").Value = "  +2.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.81"
$ws.Range("E26").Value = "  +2.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.33"
$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.938"
$ws.Range("E28").Value = "  +5.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.482"
$ws.Range("E29").Value = "  +5.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.385"
$ws.Range("E30").Value = "  +3.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09167"
$ws.Range("E31").Value = "  +1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.131"
$ws.Range("E32").Value = "  +5.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05222"
$ws.Range("E33").Value = "  +0.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7410"
$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.103"
$ws.Range("E35").Value = "  +2.42%  "

$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01834"
$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.690"
$ws.Range("E38").Value = "  +1.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9238"
$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.056"
$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4415"
$ws.Range("E41").Value = "  +2.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.939"
$ws.Range("E42").Value = "  +3.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.20"
$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.74"
$ws.Range("E44").Value = "  +22.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9939"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1386"
$ws.Range("E46").Value = "  +3.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.612"
$ws.Range("E47").Value = "  +4.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.988"
$ws.Range("E48").Value = "  +3.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.05"
$ws.Range("E49").Value = "  +5.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05828"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.428"
$ws.Range("E51").Value = "  +1.05%  "
